# Update driverAthletics Excel data for 2021-2025: add "Wins" and
# "Race Starts" columns (E, F) with data, and refresh header/row styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New header cells (E1, F1) -- existing A1:D1 headers are unchanged.
# ---------------------------------------------------------------------
$ws.Range("E1").Value2 = "Wins"
$ws.Range("F1").Value2 = "Race Starts"

# ---------------------------------------------------------------------
# 2. New data columns: Wins (E) and Race Starts (F) for rows 2-21.
#    (Name/Height/Weight/Age columns A-D are untouched.)
# ---------------------------------------------------------------------
$wins = @(103, 10, 20, 2, 53, 0, 32, 1, 8, 0, 2, 0, 1, 0, 21, 0, 0, 0, 0, 0)
$starts = @(300, 179, 140, 208, 277, 75, 330, 96, 208, 53, 63, 128, 63, 22, 349, 62, 21, 22, 59, 40)

for ($i = 0; $i -lt 20; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 5).Value2 = $wins[$i]
    $ws.Cells.Item($r, 6).Value2 = $starts[$i]
}

# ---------------------------------------------------------------------
# 3. Unified header/body styling: Arial 11, color #1B1C1D, centered both
#    ways, wrapped text -- applied across the whole used range A1:F21.
#    Build the format once on a scratch cell, then paste it across the
#    range so every cell lands on the SAME style record (avoids creating
#    one new style per previously-distinct cell format).
# ---------------------------------------------------------------------
$helper = $ws.Range("H1")
$helper.Font.Name = "Arial"
$helper.Font.Size = 11
$helper.Font.Bold = $false
$helper.Font.Color = 1907739
$helper.HorizontalAlignment = -4108
$helper.VerticalAlignment = -4108
$helper.WrapText = $true

$helper.Copy()
$ws.Range("A1:F21").PasteSpecial(-4122)
$helper.Clear()
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4. Row heights: header + most rows are 28pt, a few taller rows (wrapped
#    two-line headers/names) are 42pt.
# ---------------------------------------------------------------------
$tallRows = @(4, 10, 16, 17, 19)
for ($r = 1; $r -le 21; $r++) {
    if ($tallRows -contains $r) {
        $ws.Rows.Item($r).RowHeight = 42
    } else {
        $ws.Rows.Item($r).RowHeight = 28
    }
}

# ---------------------------------------------------------------------
# 5. View state: select F1 (also resets the scrolled-down top-left cell
#    left over from the previous edit session) and refresh the
#    used-range dimension.
# ---------------------------------------------------------------------
$ws.Range("F1").Select()

Write-Output "done"
